$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# New test-data row for the "Edit and Delete Foreign Posts" tests:
# a user whose e-mail belongs to someone else's account.
$ws.Range("A7").Value = "LoginForeinUserData"
$ws.Range("B7").Value = "aaa@abv.bg"
$ws.Range("C7").Value = "123456"

# Make the e-mail a mailto hyperlink, like the other rows in the sheet,
# then restore the shared "Hyperlink" cell style (Add() mints its own).
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:aaa@abv.bg")
$ws.Range("B7").Style = "Hyperlink"

# Reflect where the user last clicked when they saved the sheet.
$ws.Range("E11").Select()
